# The data-loader / environment-data change: column A held Excel serial
# date-time values (formatted with a custom "YYYY-MM-DD HH:MM:SS" style,
# style index 2). These are replaced with plain numeric dates encoded as
# YYYYMMDD integers, and the custom date formatting on those cells is
# removed (cells revert to the default/general format, dropping the
# custom numFmt/cellXf usage).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# (firstRow, lastRow, newYyyymmddValue) - one block per ISO week present
# in the sheet (rows 2-241).
$blocks = @(
    @(2, 17, 20180312),
    @(18, 33, 20180319),
    @(34, 49, 20180326),
    @(50, 65, 20180402),
    @(66, 81, 20180409),
    @(82, 97, 20180416),
    @(98, 113, 20180423),
    @(114, 129, 20180430),
    @(130, 145, 20180507),
    @(146, 161, 20180514),
    @(162, 177, 20180521),
    @(178, 193, 20180528),
    @(194, 209, 20180604),
    @(210, 225, 20180611),
    @(226, 241, 20180618)
)

foreach ($block in $blocks) {
    $firstRow = $block[0]
    $lastRow = $block[1]
    $newValue = $block[2]

    $rng = $ws.Range($ws.Cells.Item($firstRow, 1), $ws.Cells.Item($lastRow, 1))
    $rng.ClearFormats()
    $rng.Value = $newValue
}
